$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 19 (shifts existing rows 19:80 down to 20:81)
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 (Dia=24, total_venda=3931.75, Mes=7, Ano=2025, Periodo="07/2025")
$ws.Range("A19").Value = 24
$ws.Range("B19").Value = 3931.75
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 2025
$ws.Range("E19").Value = "07/2025"

# Corrected total_venda values for existing July rows
$ws.Range("B2").Value = 18058.83
$ws.Range("B6").Value = 24159.71
$ws.Range("B17").Value = 9021.559999999999
$ws.Range("B18").Value = 12244.93
